$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (headers) ---
# Existing header cells A1:I1 pick up the "quotePrefix" style (s="1") while
# keeping their original text. G1/H1/I1 keep their MaxDFR_* labels; a brand
# new header "PrefaultTime" is added in J1 (default style - no quote prefix).
$ws.Range("A1").Value = "'Username"
$ws.Range("B1").Value = "'Password"
$ws.Range("C1").Value = "'DeviceName"
$ws.Range("D1").Value = "'DeviceType"
$ws.Range("E1").Value = "'DeviceIPAdd"
$ws.Range("F1").Value = "'DeviceSerialNo"
$ws.Range("G1").Value = "'MaxDFR_Min"
$ws.Range("H1").Value = "'MaxDFR_Mid"
$ws.Range("I1").Value = "'MaxDFR_Max"
$ws.Range("J1").Value = "PrefaultTime"

# --- Row 2 (data) ---
# A2:E2 keep their text values, now with the quotePrefix style applied.
$ws.Range("A2").Value = "'Admin"
$ws.Range("B2").Value = "'Admin"
$ws.Range("C2").Value = "'IND_DAU_51"
$ws.Range("D2").Value = "'IDM+18"
$ws.Range("E2").Value = "'10.75.58.51"

# F2 (DeviceSerialNo value) is untouched - stays the numeric 409026540.

# New PrefaultTime value (J2), then the existing MaxDFR_* cells become
# quote-prefixed text values too - H2's underlying value changes 500 -> 1200.
$ws.Range("J2").Value = "'1000"
$ws.Range("H2").Value = "'1200"
$ws.Range("I2").Value = "'31000"
$ws.Range("G2").Value = "'400"

# --- Selection moves from G1 to H2 ---
$ws.Range("H2").Select()
